$p = $ppt.ActivePresentation

# 1. Delete the "Why Quiz app? / PART FOUR" slide (slide 8).
$p.Slides.Item(8).Delete()

# 2. The "How? / PART THREE" slide (formerly slide 10, now slide 9 after the
#    deletion above) gets its second placeholder's text changed from
#    "PART THREE" to "PART FOUR".
$s = $p.Slides.Item(9)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "PART FOUR"
